# Add two new weekly records at the top of the data table (rows 152-153),
# pushing all existing data rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 152, shifting rows 152:251 -> 154:253
$ws.Rows("152:153").Insert()

# --- New row 152 ---
$ws.Range("A152").Value = 8
$ws.Range("B152").Value = "Terminal La Palmera de La Serena"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 44767
$ws.Range("E152").Value = 4
$ws.Range("F152").Value = 100112021
$ws.Range("G152").Value = "Ají"
$ws.Range("H152").Value = "Inferno"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 600
$ws.Range("K152").Value = 14000
$ws.Range("L152").Value = 15000
$ws.Range("M152").Value = 14500
$ws.Range("N152").Value = "$/caja 12 kilos"
$ws.Range("O152").Value = "Región de Arica y Parinacota"
$ws.Range("P152").Value = 1208
$ws.Range("Q152").Value = 12
$ws.Range("R152").Value = "Hortaliza"

# --- New row 153 ---
$ws.Range("A153").Value = 8
$ws.Range("B153").Value = "Terminal La Palmera de La Serena"
$ws.Range("C153").Value = "Coquimbo"
$ws.Range("D153").Value = 44767
$ws.Range("E153").Value = 4
$ws.Range("F153").Value = 100112021
$ws.Range("G153").Value = "Ají"
$ws.Range("H153").Value = "Inferno"
$ws.Range("I153").Value = "Segunda"
$ws.Range("J153").Value = 400
$ws.Range("K153").Value = 9000
$ws.Range("L153").Value = 10000
$ws.Range("M153").Value = 9500
$ws.Range("N153").Value = "$/caja 12 kilos"
$ws.Range("O153").Value = "Región de Arica y Parinacota"
$ws.Range("P153").Value = 792
$ws.Range("Q153").Value = 12
$ws.Range("R153").Value = "Hortaliza"
